$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.05"
$ws.Range("E2").Value = "'-3.68%"
$ws.Range("E3").Value = "'-6.64%"
$ws.Range("D4").Value = "'5.106"
$ws.Range("E4").Value = "'-0.62%"
$ws.Range("D5").Value = "'0.07754"
$ws.Range("E5").Value = "'-5.85%"
$ws.Range("D6").Value = "'4.394"
$ws.Range("E6").Value = "'1.68%"
$ws.Range("D7").Value = "'1.912"
$ws.Range("E7").Value = "'-7.90%"
$ws.Range("D8").Value = "'8.225"
$ws.Range("E8").Value = "'-1.19%"
$ws.Range("E9").Value = "'-2.75%"
$ws.Range("D10").Value = "'0.9218"
$ws.Range("E10").Value = "'-2.04%"
$ws.Range("D11").Value = "'0.1268"
$ws.Range("E11").Value = "'-7.28%"
$ws.Range("D12").Value = "'0.1894"
$ws.Range("E12").Value = "'-4.45%"
$ws.Range("E13").Value = "'-3.34%"
$ws.Range("D14").Value = "'0.03447"
$ws.Range("E14").Value = "'-1.86%"
$ws.Range("D15").Value = "'0.09746"
$ws.Range("E15").Value = "'-0.73%"
$ws.Range("D16").Value = "'0.001370"
$ws.Range("E16").Value = "'-0.27%"
$ws.Range("D17").Value = "'0.006120"
$ws.Range("E17").Value = "'-2.47%"
$ws.Range("D18").Value = "'3.564"
$ws.Range("E18").Value = "'-3.56%"
$ws.Range("D19").Value = "'0.3374"
$ws.Range("E19").Value = "'-3.54%"
$ws.Range("D20").Value = "'0.1288"
$ws.Range("E20").Value = "'-1.63%"
$ws.Range("D21").Value = "'5.044"
$ws.Range("E21").Value = "'1.33%"
$ws.Range("E23").Value = "'5,218.80%"
$ws.Range("D24").Value = "'0.04345"
$ws.Range("E24").Value = "'0.03%"
$ws.Range("D25").Value = "'0.001227"
$ws.Range("E25").Value = "'0.21%"
$ws.Range("D26").Value = "'0.004494"
$ws.Range("E26").Value = "'-7.31%"
$ws.Range("D27").Value = "'0.0001360"
$ws.Range("E27").Value = "'4.89%"
$ws.Range("E39").Value = "'-1.94%"
$ws.Range("D40").Value = "'0.04939"
$ws.Range("E40").Value = "'-5.37%"
$ws.Range("D41").Value = "'0.007711"
$ws.Range("E41").Value = "'0.46%"
$ws.Range("D42").Value = "'0.009876"
$ws.Range("E42").Value = "'2.61%"
$ws.Range("D43").Value = "'0.1341"
$ws.Range("E43").Value = "'-4.76%"
$ws.Range("D44").Value = "'0.002008"
$ws.Range("E44").Value = "'-2.89%"
$ws.Range("D45").Value = "'0.008867"
$ws.Range("E45").Value = "'5.75%"
$ws.Range("D46").Value = "'0.00006841"
$ws.Range("E46").Value = "'3.59%"
$ws.Range("D47").Value = "'0.00000000756"
$ws.Range("E47").Value = "'1.22%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001310"
$ws.Range("E48").Value = "'-22.15%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003019"
$ws.Range("E49").Value = "'3.12%"
$ws.Range("D50").Value = "'0.00002116"
$ws.Range("E50").Value = "'1.22%"
$ws.Range("D51").Value = "'0.0002015"
$ws.Range("E51").Value = "'1.22%"
